{"js": "// Update the worksheet date header and all \"NN\u00d7NN=\" multiplication\n// prompts to the new values per the commit's regenerated data set.\nconst pairs = [\n  [\"2024-06-27 Thursday\", \"2024-06-28 Friday\"],\n  [\"86\u00d748=\", \"37\u00d790=\"],\n  [\"25\u00d783=\", \"47\u00d760=\"],\n  [\"86\u00d754=\", \"84\u00d790=\"],\n  [\"91\u00d766=\", \"35\u00d772=\"],\n  [\"84\u00d724=\", \"22\u00d789=\"],\n  [\"18\u00d763=\", \"42\u00d736=\"],\n  [\"81\u00d790=\", \"13\u00d723=\"],\n  [\"23\u00d756=\", \"21\u00d779=\"],\n  [\"68\u00d765=\", \"93\u00d791=\"],\n  [\"68\u00d741=\", \"74\u00d726=\"],\n  [\"88\u00d793=\", \"86\u00d791=\"],\n  [\"80\u00d786=\", \"61\u00d798=\"],\n  [\"53\u00d741=\", \"36\u00d797=\"],\n  [\"51\u00d758=\", \"19\u00d794=\"],\n  [\"26\u00d783=\", \"75\u00d775=\"],\n  [\"75\u00d780=\", \"89\u00d797=\"],\n  [\"47\u00d797=\", \"78\u00d739=\"],\n  [\"28\u00d733=\", \"74\u00d752=\"],\n  [\"29\u00d777=\", \"57\u00d716=\"],\n  [\"32\u00d728=\", \"59\u00d755=\"],\n  [\"51\u00d711=\", \"48\u00d716=\"],\n  [\"83\u00d749=\", \"19\u00d774=\"],\n  [\"20\u00d774=\", \"19\u00d771=\"],\n  [\"93\u00d790=\", \"87\u00d768=\"],\n  [\"96\u00d798=\", \"81\u00d742=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and all \"NN\u00d7NN=\" multiplication\n# prompts to the new values per the commit's regenerated data set.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-06-27 Thursday\", \"2024-06-28 Friday\"),\n    @(\"86\u00d748=\", \"37\u00d790=\"),\n    @(\"25\u00d783=\", \"47\u00d760=\"),\n    @(\"86\u00d754=\", \"84\u00d790=\"),\n    @(\"91\u00d766=\", \"35\u00d772=\"),\n    @(\"84\u00d724=\", \"22\u00d789=\"),\n    @(\"18\u00d763=\", \"42\u00d736=\"),\n    @(\"81\u00d790=\", \"13\u00d723=\"),\n    @(\"23\u00d756=\", \"21\u00d779=\"),\n    @(\"68\u00d765=\", \"93\u00d791=\"),\n    @(\"68\u00d741=\", \"74\u00d726=\"),\n    @(\"88\u00d793=\", \"86\u00d791=\"),\n    @(\"80\u00d786=\", \"61\u00d798=\"),\n    @(\"53\u00d741=\", \"36\u00d797=\"),\n    @(\"51\u00d758=\", \"19\u00d794=\"),\n    @(\"26\u00d783=\", \"75\u00d775=\"),\n    @(\"75\u00d780=\", \"89\u00d797=\"),\n    @(\"47\u00d797=\", \"78\u00d739=\"),\n    @(\"28\u00d733=\", \"74\u00d752=\"),\n    @(\"29\u00d777=\", \"57\u00d716=\"),\n    @(\"32\u00d728=\", \"59\u00d755=\"),\n    @(\"51\u00d711=\", \"48\u00d716=\"),\n    @(\"83\u00d749=\", \"19\u00d774=\"),\n    @(\"20\u00d774=\", \"19\u00d771=\"),\n    @(\"93\u00d790=\", \"87\u00d768=\"),\n    @(\"96\u00d798=\", \"81\u00d742=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap = wdFindContinue\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace = wdReplaceAll\n    )\n}\n"}
